# DSM.xlsx update: add report/flex-table configuration rows to the
# "configurations" sheet, and switch the active sheet/selection back to
# "configurations" (from "decision_makers_options").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "configurations"

# Scratch cell used to coerce the literal text "True"/"False" into the
# sheet as a genuine *text* value (t="s") instead of letting a plain
# Value assignment of "True"/"False" get auto-coerced into a boolean
# (t="b") by the engine. Building it via a quoted-string formula and
# then Copy + PasteSpecial(values) drops the formula and yields a plain
# shared-string text cell with no extra "quote prefix" styling.
$scratch = $ws.Cells.Item(100, 10)

function Set-TextBool($cell, [string]$val) {
    $scratch.Formula = "=""" + $val + """"
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
}

$ws.Cells.Item(3, 1).Value = "Optimize_DMO_name"

$ws.Cells.Item(4, 1).Value = "report_title_page"
Set-TextBool $ws.Cells.Item(4, 2) "True"

$ws.Cells.Item(5, 1).Value = "report_strategic_challenge"
Set-TextBool $ws.Cells.Item(5, 2) "True"

$ws.Cells.Item(6, 1).Value = "report_key_outputs_theme"
Set-TextBool $ws.Cells.Item(6, 2) "True"

$ws.Cells.Item(7, 1).Value = "report_decision_makers_options"
Set-TextBool $ws.Cells.Item(7, 2) "True"

$ws.Cells.Item(8, 1).Value = "report_scenarios"
Set-TextBool $ws.Cells.Item(8, 2) "True"

$ws.Cells.Item(9, 1).Value = "report_fixed_inputs"
Set-TextBool $ws.Cells.Item(9, 2) "True"

$ws.Cells.Item(10, 1).Value = "report_dependencies"
Set-TextBool $ws.Cells.Item(10, 2) "False"

$ws.Cells.Item(11, 1).Value = "report_weighted_appreciations"
Set-TextBool $ws.Cells.Item(11, 2) "True"

$ws.Cells.Item(12, 1).Value = "report_add_optimize"
Set-TextBool $ws.Cells.Item(12, 2) "False"

# Remove the scratch helper cell again.
$scratch.Clear()

# Move the active tab / selection from "decision_makers_options" back to
# "configurations" at cell C17.
$ws.Select() | Out-Null
$ws.Range("C17").Select() | Out-Null

Write-Output "DSM.xlsx configurations sheet updated"
